# ---------------------------------------------------------------------------
# Applies the resume content refresh described in the target diff:
#   1. Collapses the three detailed CORE COMPETENCIES paragraphs into a
#      single one-line summary.
#   2. Replaces the bullet lists under five PROFESSIONAL EXPERIENCE entries
#      (RESEARCH DIRECTOR, SOFTWARE ENGINEER - Salsa Labs, INTERIM
#      TECHNOLOGY MANAGER, PROGRAMMER - Lake Research Partners, FIELD
#      DIRECTOR) with new, more specific bullet text (several sections gain
#      extra bullets).
#   3. Appends a new "TECHNICAL SKILLS" section (Heading2 + three summary
#      paragraphs) at the very end of the document.
# ---------------------------------------------------------------------------

function AssertParaText($d, $index, $expectedSubstring) {
    $actual = $d.Paragraphs.Item($index).Range.Text
    if ($actual.IndexOf($expectedSubstring) -lt 0) {
        throw "Paragraph $index did not contain expected text '$expectedSubstring'. Actual: '$actual'"
    }
}

function SetText($para, $text) {
    $para.Range.Text = $text
}

# Insert a new paragraph containing $text immediately after $para (which
# keeps $para's own paragraph mark/style intact) and returns nothing -
# the new paragraph picks up $para's style since it's created by splitting
# $para's trailing paragraph mark.
function InsertParaAfter($d, $para, $text) {
    $r = $d.Range($para.Range.Start, $para.Range.End - 1)
    $r.Collapse(0)
    $r.InsertAfter([char]13 + $text)
}

# Replace the block of $oldCount existing paragraphs starting at
# $startIndex with the strings in $newTexts (array). Handles both growing
# and shrinking the number of paragraphs in the block.
function ReplaceBlock($d, $startIndex, $oldCount, $newTexts) {
    $n = $newTexts.Count
    $reuse = [Math]::Min($oldCount, $n)

    for ($i = 0; $i -lt $reuse; $i++) {
        $p = $d.Paragraphs.Item($startIndex + $i)
        SetText $p $newTexts[$i]
    }

    if ($oldCount -gt $n) {
        $toDelete = $oldCount - $n
        for ($i = 0; $i -lt $toDelete; $i++) {
            $d.Paragraphs.Item($startIndex + $reuse).Range.Delete()
        }
    } elseif ($n -gt $oldCount) {
        $lastPara = $d.Paragraphs.Item($startIndex + $reuse - 1)
        for ($i = $reuse; $i -lt $n; $i++) {
            InsertParaAfter $d $lastPara $newTexts[$i]
            $lastPara = $d.Paragraphs.Item($startIndex + $i)
        }
    }
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. CORE COMPETENCIES: collapse the three detailed bullet paragraphs
#    (indices 6-8) into one short summary paragraph.
# ---------------------------------------------------------------------------
AssertParaText $d 6 "Research and Analytics: Survey Methodology"
AssertParaText $d 7 "Programming and Development: Python"
AssertParaText $d 8 "Data Infrastructure: Cloud Platforms"

SetText $d.Paragraphs.Item(6) "Research and Analytics • Programming and Development • Data Infrastructure"
$d.Paragraphs.Item(7).Range.Delete()
$d.Paragraphs.Item(7).Range.Delete()

# ---------------------------------------------------------------------------
# 2. RESEARCH DIRECTOR - Progressive Change Campaign Committee: 4 -> 6
#    bullets, all new text. (indices shifted -2 by step 1, so 39-42 -> 37-40)
# ---------------------------------------------------------------------------
AssertParaText $d 37 "Managed critical research operations"
$rd = @(
  "• Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls",
  "• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren",
  "• Built tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver",
  "• Designed survey deployment system facilitating thousands of simultaneous phone surveys",
  "• Significantly increased data collection efficiency through automated calling infrastructure",
  "• Managed comprehensive research operations for progressive political initiatives and candidates"
)
ReplaceBlock $d 37 4 $rd

# ---------------------------------------------------------------------------
# 3. SOFTWARE ENGINEER - Salsa Labs: 4 -> 5 bullets, all new text.
# ---------------------------------------------------------------------------
AssertParaText $d 45 "Developed software solutions for political campaigns"
$salsa = @(
  "• Maintained and extended entire geospatial analysis and reporting tools for Java-based CRM system",
  "• Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers",
  "• Built geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill",
  "• Integrated mapping and visualization tools for political campaign data analysis",
  "• Collaborated with political strategists to translate geospatial requirements into technical solutions"
)
ReplaceBlock $d 45 4 $salsa

# ---------------------------------------------------------------------------
# 4. INTERIM TECHNOLOGY MANAGER - The Praxis Project: 4 -> 6 bullets.
# ---------------------------------------------------------------------------
AssertParaText $d 52 "Integrated technology solutions within organizational frameworks"
$praxis = @(
  "• Assisted in search for full-time CTO while performing all programmatic technology roles for multi-million dollar organization",
  "• Made all technology decisions and practices for massive multinational non-governmental organization",
  "• Wrote comprehensive frameworks for internal and external technology audits",
  "• Trained beneficiaries on spatial and Census data analysis for public health research",
  "• Trained NGO staff in web development using Drupal, PHP, and MySQL",
  "• Managed technology infrastructure supporting community health initiatives across multiple countries"
)
ReplaceBlock $d 52 4 $praxis

# ---------------------------------------------------------------------------
# 5. PROGRAMMER - Lake Research Partners: 4 -> 6 bullets.
# ---------------------------------------------------------------------------
AssertParaText $d 60 "Developed data analysis tools for political polling"
$lake = @(
  "• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party",
  "• Developed system that later became the Polling Consortium Database at The Analyst Institute",
  "• Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections",
  "• Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle",
  "• Pioneered integration of advanced mapping techniques into standard reports including choropleths and hexagonal grid maps",
  "• Developed innovative approaches to visualizing demographic and market data for enhanced client understanding"
)
ReplaceBlock $d 60 4 $lake

# ---------------------------------------------------------------------------
# 6. FIELD DIRECTOR - The Feldman Group: 4 -> 6 bullets.
# ---------------------------------------------------------------------------
AssertParaText $d 68 "Managed field operations for political campaigns"
$feldman = @(
  "• Administered all quantitative and qualitative research operations ensuring reporting accuracy",
  "• Managed comprehensive survey fielding for multi-million dollar research firm",
  "• Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings",
  "• Created custom reports and data visualizations based on specific client requirements",
  "• Introduced mapping and geospatial analysis into standard reporting procedures",
  "• Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL"
)
ReplaceBlock $d 68 4 $feldman

# ---------------------------------------------------------------------------
# 7. Append a new TECHNICAL SKILLS section at the end of the document.
#    All four new paragraphs are inserted in a single InsertAfter call so
#    they inherit the "Normal" style of the final existing bullet
#    paragraph (avoiding an explicit pStyle on the three body paragraphs);
#    only the first of the four is then promoted to Heading2.
# ---------------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
AssertParaText $d $lastIndex "Redistricting analysis used in court cases"

$lastPara = $d.Paragraphs.Item($lastIndex)
$r = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$r.Collapse(0)
$newBlock = "TECHNICAL SKILLS" + [char]13 + `
    "RESEARCH AND ANALYTICS Survey Methodology; Statistical Analysis; Geospatial Analysis; Data Visualization; Research Management" + [char]13 + `
    "PROGRAMMING AND DEVELOPMENT Python; JVM Languages; Web Technologies; Database Languages; Statistical Computing" + [char]13 + `
    "DATA INFRASTRUCTURE Cloud Platforms; Big Data; Databases; Geospatial; DevOps"
$r.InsertAfter([char]13 + $newBlock)

$headingPara = $d.Paragraphs.Item($lastIndex + 1)
$headingPara.Range.ParagraphFormat.set_Style("Heading 2")

Write-Host "Done. Total paragraphs: " $d.Paragraphs.Count
